$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 22675.861
$ws.Range("J87").Value = 22675.861
$ws.Range("L87").Value = 22675.861
$ws.Range("N87").Value = -25171.861

$ws.Range("H90").Value = 22675.861
$ws.Range("J90").Value = 22675.861
$ws.Range("L90").Value = 68027.583
$ws.Range("N90").Value = -80507.583

$ws.Range("H99").Value = 1709.8572
$ws.Range("I99").Value = 272.66666
$ws.Range("J99").Value = 2787.75
$ws.Range("K99").Value = 817.9999799999999
$ws.Range("L99").Value = 8363.25
$ws.Range("M99").Value = 680.0000200000001
$ws.Range("N99").Value = -11359.25

$ws.Range("H106").Value = 2390.9092
$ws.Range("I106").Value = 580
$ws.Range("J106").Value = 3900
$ws.Range("K106").Value = 580
$ws.Range("L106").Value = 3900
$ws.Range("M106").Value = 51
$ws.Range("N106").Value = -5162

$ws.Range("H127").Value = 616.13043
$ws.Range("I127").Value = 616.6667
$ws.Range("J127").Value = 616.05
$ws.Range("K127").Value = 1850.0001
$ws.Range("L127").Value = 1848.15
$ws.Range("M127").Value = 3109.9999
$ws.Range("N127").Value = -11768.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6595.654
$ws.Range("I32").Value = 5554.3286
$ws.Range("J32").Value = 21799
$ws.Range("K32").Value = 5554.3286
$ws.Range("L32").Value = 21799
$ws.Range("M32").Value = -5267.3286
$ws.Range("N32").Value = -22373

$ws.Range("H74").Value = 1160.8889
$ws.Range("I74").Value = 1192.1333
$ws.Range("J74").Value = 1004.6667
$ws.Range("K74").Value = 1192.1333
$ws.Range("L74").Value = 1004.6667
$ws.Range("M74").Value = -318.1333
$ws.Range("N74").Value = -2752.6667

$ws.Range("H77").Value = 1160.8889
$ws.Range("I77").Value = 1192.1333
$ws.Range("J77").Value = 1004.6667
$ws.Range("K77").Value = 5960.666499999999
$ws.Range("L77").Value = 5023.3335
$ws.Range("M77").Value = -1592.666499999999
$ws.Range("N77").Value = -13759.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 696.55554
$ws.Range("I94").Value = 719.875
$ws.Range("J94").Value = 510
$ws.Range("K94").Value = 719.875
$ws.Range("L94").Value = 510
$ws.Range("M94").Value = -268.875
$ws.Range("N94").Value = -1412

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 65844050
$ws.Range("I5").Value = 197531150
$ws.Range("K5").Value = 197531150
$ws.Range("M5").Value = -197531038

$ws.Range("H8").Value = 755
$ws.Range("J8").Value = 1010
$ws.Range("L8").Value = 1010
$ws.Range("N8").Value = -1290

$ws.Range("H16").Value = 2036.3636
$ws.Range("I16").Value = 1180
$ws.Range("J16").Value = 2750
$ws.Range("K16").Value = 1180
$ws.Range("L16").Value = 2750
$ws.Range("M16").Value = -893
$ws.Range("N16").Value = -3324

$ws.Range("H58").Value = 14709539
$ws.Range("I58").Value = 2517.0417
$ws.Range("J58").Value = 50006390
$ws.Range("K58").Value = 2517.0417
$ws.Range("L58").Value = 50006390
$ws.Range("M58").Value = -2314.0417
$ws.Range("N58").Value = -50006796

$ws.Range("H107").Value = 1405.8636
$ws.Range("I107").Value = 516.3125
$ws.Range("K107").Value = 516.3125
$ws.Range("M107").Value = 1403.6875

$ws.Range("H113").Value = 2036.3636
$ws.Range("I113").Value = 1180
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 1180
$ws.Range("L113").Value = 2750
$ws.Range("M113").Value = 990
$ws.Range("N113").Value = -7090

$ws.Range("H122").Value = 2534.3333
$ws.Range("I122").Value = 1948.6522
$ws.Range("J122").Value = 3881.4
$ws.Range("K122").Value = 5845.9566
$ws.Range("L122").Value = 11644.2
$ws.Range("M122").Value = -3395.9566
$ws.Range("N122").Value = -16544.2

$ws.Range("H132").Value = 2150.3
$ws.Range("I132").Value = 1653.7333
$ws.Range("J132").Value = 3640
$ws.Range("K132").Value = 4961.199900000001
$ws.Range("L132").Value = 10920
$ws.Range("M132").Value = -2431.199900000001
$ws.Range("N132").Value = -15980

$ws.Range("H136").Value = 14709539
$ws.Range("I136").Value = 2517.0417
$ws.Range("J136").Value = 50006390
$ws.Range("K136").Value = 7551.125100000001
$ws.Range("L136").Value = 150019170
$ws.Range("M136").Value = -5001.125100000001
$ws.Range("N136").Value = -150024270

$ws.Range("H138").Value = 27447.47
$ws.Range("J138").Value = 27447.47
$ws.Range("L138").Value = 27447.47
$ws.Range("N138").Value = -37727.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1809.8889
$ws.Range("J55").Value = 2757.8
$ws.Range("L55").Value = 8273.400000000001
$ws.Range("N55").Value = -8627.400000000001

$ws.Range("H87").Value = 9997
$ws.Range("I87").Value = 1993
$ws.Range("K87").Value = 5979
$ws.Range("M87").Value = -4731

$ws.Range("H90").Value = 9997
$ws.Range("I90").Value = 1993
$ws.Range("K90").Value = 17937
$ws.Range("M90").Value = -11697

$ws.Range("H129").Value = 2016
$ws.Range("I129").Value = 1594.2858
$ws.Range("J129").Value = 3000
$ws.Range("K129").Value = 4782.857400000001
$ws.Range("L129").Value = 9000
$ws.Range("M129").Value = 217.1425999999992
$ws.Range("N129").Value = -19000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 35289.582
$ws.Range("I102").Value = 2404.6843
$ws.Range("J102").Value = 87357.336
$ws.Range("K102").Value = 2404.6843
$ws.Range("L102").Value = 87357.336
$ws.Range("M102").Value = -782.6842999999999
$ws.Range("N102").Value = -90601.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3126280.8
$ws.Range("I7").Value = 4762682
$ws.Range("K7").Value = 4762682
$ws.Range("M7").Value = -4762570

$ws.Range("H93").Value = 1695.1818
$ws.Range("I93").Value = 830.375
$ws.Range("K93").Value = 830.375
$ws.Range("M93").Value = 417.625

$ws.Range("H126").Value = 3126280.8
$ws.Range("I126").Value = 4762682
$ws.Range("K126").Value = 14288046
$ws.Range("M126").Value = -14285576

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 314718.8
$ws.Range("I122").Value = 456233.28
$ws.Range("K122").Value = 1368699.84
$ws.Range("M122").Value = -1366249.84

$ws.Range("H126").Value = 3032061.5
$ws.Range("I126").Value = 1021.75
$ws.Range("K126").Value = 3065.25
$ws.Range("M126").Value = -595.25

Write-Output "done"